# Adds a new data row (row 31) to the "Artfynd" sheet, matching the
# source record for "Höstlåsbräken" (Botrychium multifidum) found at
# Bödagården NO, Öl on 2023-09-03.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain numeric cells -------------------------------------------------
$ws.Range("A31").Value = 111885493
$ws.Range("B31").Value = 95610
$ws.Range("E31").Value = 167
$ws.Range("Q31").Value = 623795.9600775555
$ws.Range("R31").Value = 6348423.740576888
$ws.Range("S31").Value = 25

# --- Plain text cells ------------------------------------------------------
$ws.Range("C31").Value = "Ovaliderad"
$ws.Range("D31").Value = "NT"
$ws.Range("F31").Value = "Höstlåsbräken"
$ws.Range("G31").Value = "Botrychium multifidum"
$ws.Range("H31").Value = "(S. G. Gmel.) Rupr."
$ws.Range("P31").Value = "Bödagården NO, Öl"
$ws.Range("T31").Value = "Kalmar"
$ws.Range("U31").Value = "Borgholm"
$ws.Range("V31").Value = "Öland"
$ws.Range("W31").Value = "Böda"
$ws.Range("Z31").Value = "00:00"
$ws.Range("AB31").Value = "00:00"
$ws.Range("AI31").Value = "nyröjt område mellan sanddyner"
$ws.Range("AW31").Value = "Ulla-Britt Andersson"
$ws.Range("AX31").Value = "Ulla-Britt Andersson, Thomas Gunnarsson"

# --- Date-looking cells that must stay plain text, not auto-converted to
# Excel date serials: force Text formatting for the assignment, then
# restore the default "Normal" style so no leftover number format sticks
# to the cell (keeps the cell a plain text cell, styled like the rest of
# the sheet).
$ws.Range("Y31").NumberFormat = "@"
$ws.Range("Y31").Value = "2023-09-03"
$ws.Range("Y31").Style = "Normal"

$ws.Range("AA31").NumberFormat = "@"
$ws.Range("AA31").Value = "2023-09-03"
$ws.Range("AA31").Style = "Normal"

# --- Boolean cells -----------------------------------------------------
$ws.Range("AD31").Value = $true
$ws.Range("AE31").Value = $false
$ws.Range("AG31").Value = $false
